$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before EF (shifts EF:FJ -> EG:FK), matching the
# workbook's dimension growing from A1:FJ25 to A1:FK25.
$ws.Range("EF1").EntireColumn.Insert()

# Header row: new column gets the next day in the "-dec" sequence.
$ws.Range("EF1").Value = "09-dec"

# Data rows: new column is filled with the "no data" placeholder used
# throughout the rest of the sheet.
$ws.Range("EF2:EF25").Value = "-"
